# Auto-generated edit script for csv/maps.xlsx update
# Applies: dimension growth, several cell updates in existing rows,
# and 17 new data rows (1456-1472) appended to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $r, $c, $val) {
    # Force the cell to store $val as literal text (avoids Excel
    # auto-converting date-like strings, e.g. "2021/8/18", into date serials).
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

function Set-PlainValue($ws, $r, $c, $val) {
    $ws.Cells.Item($r, $c).Value = $val
}

function Set-EmptyTextCell($ws, $r, $c) {
    # Create an (empty) placeholder cell, matching the source data
    # convention where every column always has a cell, even if blank.
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = ""
    $cell.Style = "Normal"
}

# ---- Updates to existing rows ----
Set-TextValue $ws 1316 11 "2021/8/17"
Set-PlainValue $ws 1326 12 60
Set-PlainValue $ws 1342 12 35
Set-PlainValue $ws 1376 12 134
Set-TextValue $ws 1377 11 "2021/8/17"
Set-PlainValue $ws 1386 12 12
Set-TextValue $ws 1390 11 "2021/8/17"
Set-PlainValue $ws 1394 12 13
Set-PlainValue $ws 1431 12 3
Set-PlainValue $ws 1436 12 17
Set-PlainValue $ws 1440 12 15
Set-PlainValue $ws 1442 12 2
Set-PlainValue $ws 1447 2 43.02468969
Set-PlainValue $ws 1447 3 141.3468009
Set-PlainValue $ws 1447 8 "陸上自衛隊札幌駐屯地(市384例)"
Set-PlainValue $ws 1450 12 14

# ---- New rows 1456-1472 ----
# row 1456
Set-PlainValue $ws 1456 1 1455
Set-PlainValue $ws 1456 2 41.82240883
Set-PlainValue $ws 1456 3 140.6394737
Set-PlainValue $ws 1456 4 "北海道"
Set-PlainValue $ws 1456 5 "渡島地方"
Set-PlainValue $ws 1456 6 "事業所"
Set-PlainValue $ws 1456 7 3
Set-PlainValue $ws 1456 8 "遊技施設"
Set-PlainValue $ws 1456 9 1
Set-TextValue $ws 1456 10 "2021/8/18"
Set-EmptyTextCell $ws 1456 11
Set-PlainValue $ws 1456 12 7
Set-PlainValue $ws 1456 13 0
Set-PlainValue $ws 1456 14 1
Set-PlainValue $ws 1456 15 0
Set-EmptyTextCell $ws 1456 16
Set-PlainValue $ws 1456 17 "共用室"
Set-PlainValue $ws 1456 18 "その他"
Set-PlainValue $ws 1456 19 3
# row 1457
Set-PlainValue $ws 1457 1 1456
Set-PlainValue $ws 1457 2 42.76334649
Set-PlainValue $ws 1457 3 143.0311555
Set-PlainValue $ws 1457 4 "北海道"
Set-PlainValue $ws 1457 5 "帯広市"
Set-PlainValue $ws 1457 6 "保育園/幼稚園"
Set-PlainValue $ws 1457 7 16
Set-PlainValue $ws 1457 8 "あけぼの保育園"
Set-PlainValue $ws 1457 9 1
Set-TextValue $ws 1457 10 "2021/8/18"
Set-EmptyTextCell $ws 1457 11
Set-PlainValue $ws 1457 12 23
Set-PlainValue $ws 1457 13 0
Set-PlainValue $ws 1457 14 1
Set-PlainValue $ws 1457 15 0
Set-EmptyTextCell $ws 1457 16
Set-PlainValue $ws 1457 17 "共用室"
Set-PlainValue $ws 1457 18 "教育/保育"
Set-PlainValue $ws 1457 19 16
# row 1458
Set-PlainValue $ws 1458 1 1457
Set-PlainValue $ws 1458 2 43.04297113
Set-PlainValue $ws 1458 3 144.3165644
Set-PlainValue $ws 1458 4 "北海道"
Set-PlainValue $ws 1458 5 "釧路地方"
Set-PlainValue $ws 1458 6 "事業所"
Set-PlainValue $ws 1458 7 3
Set-PlainValue $ws 1458 8 "食品加工業の事業所"
Set-PlainValue $ws 1458 9 1
Set-TextValue $ws 1458 10 "2021/8/18"
Set-EmptyTextCell $ws 1458 11
Set-PlainValue $ws 1458 12 7
Set-PlainValue $ws 1458 13 0
Set-PlainValue $ws 1458 14 1
Set-PlainValue $ws 1458 15 0
Set-EmptyTextCell $ws 1458 16
Set-PlainValue $ws 1458 17 "食品加工場"
Set-PlainValue $ws 1458 18 "その他"
Set-PlainValue $ws 1458 19 19
# row 1459
Set-PlainValue $ws 1459 1 1458
Set-PlainValue $ws 1459 2 43.03294396
Set-PlainValue $ws 1459 3 141.3128175
Set-PlainValue $ws 1459 4 "北海道"
Set-PlainValue $ws 1459 5 "札幌市"
Set-PlainValue $ws 1459 6 "事業所"
Set-PlainValue $ws 1459 7 3
Set-PlainValue $ws 1459 8 "コールセンター(市385例)"
Set-PlainValue $ws 1459 9 1
Set-TextValue $ws 1459 10 "2021/8/18"
Set-EmptyTextCell $ws 1459 11
Set-PlainValue $ws 1459 12 5
Set-PlainValue $ws 1459 13 0
Set-PlainValue $ws 1459 14 1
Set-PlainValue $ws 1459 15 0
Set-EmptyTextCell $ws 1459 16
Set-PlainValue $ws 1459 17 "コールセンター"
Set-PlainValue $ws 1459 18 "その他"
Set-PlainValue $ws 1459 19 23
# row 1460
Set-PlainValue $ws 1460 1 1459
Set-PlainValue $ws 1460 2 43.05652926
Set-PlainValue $ws 1460 3 141.2860383
Set-PlainValue $ws 1460 4 "北海道"
Set-PlainValue $ws 1460 5 "札幌市"
Set-PlainValue $ws 1460 6 "保育園/幼稚園"
Set-PlainValue $ws 1460 7 16
Set-PlainValue $ws 1460 8 "認可保育施設(市386例)"
Set-PlainValue $ws 1460 9 1
Set-TextValue $ws 1460 10 "2021/8/18"
Set-EmptyTextCell $ws 1460 11
Set-PlainValue $ws 1460 12 5
Set-PlainValue $ws 1460 13 0
Set-PlainValue $ws 1460 14 1
Set-PlainValue $ws 1460 15 0
Set-EmptyTextCell $ws 1460 16
Set-PlainValue $ws 1460 17 "共用室"
Set-PlainValue $ws 1460 18 "教育/保育"
Set-PlainValue $ws 1460 19 16
# row 1461
Set-PlainValue $ws 1461 1 1460
Set-PlainValue $ws 1461 2 43.06241095
Set-PlainValue $ws 1461 3 141.3543572
Set-PlainValue $ws 1461 4 "北海道"
Set-PlainValue $ws 1461 5 "札幌市"
Set-PlainValue $ws 1461 6 "事業所"
Set-PlainValue $ws 1461 7 3
Set-PlainValue $ws 1461 8 "環境局 環境事業部(家庭ごみ収拾運搬)"
Set-PlainValue $ws 1461 9 1
Set-TextValue $ws 1461 10 "2021/8/18"
Set-EmptyTextCell $ws 1461 11
Set-PlainValue $ws 1461 12 1
Set-PlainValue $ws 1461 13 0
Set-PlainValue $ws 1461 14 1
Set-PlainValue $ws 1461 15 0
Set-EmptyTextCell $ws 1461 16
Set-PlainValue $ws 1461 17 "オフィス"
Set-PlainValue $ws 1461 18 "その他"
Set-PlainValue $ws 1461 19 3
# row 1462
Set-PlainValue $ws 1462 1 1461
Set-PlainValue $ws 1462 2 43.11327184
Set-PlainValue $ws 1462 3 141.3670795
Set-PlainValue $ws 1462 4 "北海道"
Set-PlainValue $ws 1462 5 "札幌市"
Set-PlainValue $ws 1462 6 "事業所"
Set-PlainValue $ws 1462 7 3
Set-PlainValue $ws 1462 8 "地下鉄栄町駅〜つどーむ間 乗合タクシー運転手"
Set-PlainValue $ws 1462 9 1
Set-TextValue $ws 1462 10 "2021/8/18"
Set-EmptyTextCell $ws 1462 11
Set-PlainValue $ws 1462 12 1
Set-PlainValue $ws 1462 13 0
Set-PlainValue $ws 1462 14 1
Set-PlainValue $ws 1462 15 0
Set-EmptyTextCell $ws 1462 16
Set-PlainValue $ws 1462 17 "自家用車"
Set-PlainValue $ws 1462 18 "その他"
Set-PlainValue $ws 1462 19 24
# row 1463
Set-PlainValue $ws 1463 1 1462
Set-PlainValue $ws 1463 2 43.05706823
Set-PlainValue $ws 1463 3 141.3804845
Set-PlainValue $ws 1463 4 "北海道"
Set-PlainValue $ws 1463 5 "札幌市"
Set-PlainValue $ws 1463 6 "中学校"
Set-PlainValue $ws 1463 7 7
Set-PlainValue $ws 1463 8 "札幌市立幌東中学校"
Set-PlainValue $ws 1463 9 1
Set-TextValue $ws 1463 10 "2021/8/18"
Set-EmptyTextCell $ws 1463 11
Set-PlainValue $ws 1463 12 1
Set-PlainValue $ws 1463 13 0
Set-PlainValue $ws 1463 14 1
Set-PlainValue $ws 1463 15 0
Set-EmptyTextCell $ws 1463 16
Set-PlainValue $ws 1463 17 "その他"
Set-PlainValue $ws 1463 18 "教育/保育"
Set-PlainValue $ws 1463 19 7
# row 1464
Set-PlainValue $ws 1464 1 1463
Set-PlainValue $ws 1464 2 43.04260464
Set-PlainValue $ws 1464 3 141.3998173
Set-PlainValue $ws 1464 4 "北海道"
Set-PlainValue $ws 1464 5 "札幌市"
Set-PlainValue $ws 1464 6 "小学校"
Set-PlainValue $ws 1464 7 6
Set-PlainValue $ws 1464 8 "札幌市立南白石小学校"
Set-PlainValue $ws 1464 9 1
Set-TextValue $ws 1464 10 "2021/8/18"
Set-EmptyTextCell $ws 1464 11
Set-PlainValue $ws 1464 12 1
Set-PlainValue $ws 1464 13 0
Set-PlainValue $ws 1464 14 1
Set-PlainValue $ws 1464 15 0
Set-EmptyTextCell $ws 1464 16
Set-PlainValue $ws 1464 17 "その他"
Set-PlainValue $ws 1464 18 "教育/保育"
Set-PlainValue $ws 1464 19 6
# row 1465
Set-PlainValue $ws 1465 1 1464
Set-PlainValue $ws 1465 2 42.94663432
Set-PlainValue $ws 1465 3 141.3480175
Set-PlainValue $ws 1465 4 "北海道"
Set-PlainValue $ws 1465 5 "札幌市"
Set-PlainValue $ws 1465 6 "中学校"
Set-PlainValue $ws 1465 7 7
Set-PlainValue $ws 1465 8 "札幌市立常盤中学校"
Set-PlainValue $ws 1465 9 1
Set-TextValue $ws 1465 10 "2021/8/18"
Set-EmptyTextCell $ws 1465 11
Set-PlainValue $ws 1465 12 1
Set-PlainValue $ws 1465 13 0
Set-PlainValue $ws 1465 14 1
Set-PlainValue $ws 1465 15 0
Set-EmptyTextCell $ws 1465 16
Set-PlainValue $ws 1465 17 "その他"
Set-PlainValue $ws 1465 18 "教育/保育"
Set-PlainValue $ws 1465 19 7
# row 1466
Set-PlainValue $ws 1466 1 1465
Set-PlainValue $ws 1466 2 43.07124885
Set-PlainValue $ws 1466 3 141.4193157
Set-PlainValue $ws 1466 4 "北海道"
Set-PlainValue $ws 1466 5 "札幌市"
Set-PlainValue $ws 1466 6 "中学校"
Set-PlainValue $ws 1466 7 7
Set-PlainValue $ws 1466 8 "札幌市立米里中学校"
Set-PlainValue $ws 1466 9 1
Set-TextValue $ws 1466 10 "2021/8/18"
Set-EmptyTextCell $ws 1466 11
Set-PlainValue $ws 1466 12 1
Set-PlainValue $ws 1466 13 0
Set-PlainValue $ws 1466 14 1
Set-PlainValue $ws 1466 15 0
Set-EmptyTextCell $ws 1466 16
Set-PlainValue $ws 1466 17 "その他"
Set-PlainValue $ws 1466 18 "教育/保育"
Set-PlainValue $ws 1466 19 7
# row 1467
Set-PlainValue $ws 1467 1 1466
Set-PlainValue $ws 1467 2 43.04858557
Set-PlainValue $ws 1467 3 141.3723059
Set-PlainValue $ws 1467 4 "北海道"
Set-PlainValue $ws 1467 5 "札幌市"
Set-PlainValue $ws 1467 6 "小学校"
Set-PlainValue $ws 1467 7 6
Set-PlainValue $ws 1467 8 "札幌市立豊平小学校"
Set-PlainValue $ws 1467 9 1
Set-TextValue $ws 1467 10 "2021/8/18"
Set-EmptyTextCell $ws 1467 11
Set-PlainValue $ws 1467 12 1
Set-PlainValue $ws 1467 13 0
Set-PlainValue $ws 1467 14 1
Set-PlainValue $ws 1467 15 0
Set-EmptyTextCell $ws 1467 16
Set-PlainValue $ws 1467 17 "その他"
Set-PlainValue $ws 1467 18 "教育/保育"
Set-PlainValue $ws 1467 19 6
# row 1468
Set-PlainValue $ws 1468 1 1467
Set-PlainValue $ws 1468 2 43.81090131
Set-PlainValue $ws 1468 3 142.4368403
Set-PlainValue $ws 1468 4 "北海道"
Set-PlainValue $ws 1468 5 "旭川市"
Set-PlainValue $ws 1468 6 "事業所"
Set-PlainValue $ws 1468 7 3
Set-PlainValue $ws 1468 8 "旭川市 永山支所"
Set-PlainValue $ws 1468 9 1
Set-TextValue $ws 1468 10 "2021/8/18"
Set-EmptyTextCell $ws 1468 11
Set-PlainValue $ws 1468 12 1
Set-PlainValue $ws 1468 13 0
Set-PlainValue $ws 1468 14 1
Set-PlainValue $ws 1468 15 0
Set-EmptyTextCell $ws 1468 16
Set-PlainValue $ws 1468 17 "オフィス"
Set-PlainValue $ws 1468 18 "その他"
Set-PlainValue $ws 1468 19 3
# row 1469
Set-PlainValue $ws 1469 1 1468
Set-PlainValue $ws 1469 2 42.33895451
Set-PlainValue $ws 1469 3 142.3750741
Set-PlainValue $ws 1469 4 "北海道"
Set-PlainValue $ws 1469 5 "新ひだか町"
Set-PlainValue $ws 1469 6 "医療機関"
Set-PlainValue $ws 1469 7 5
Set-PlainValue $ws 1469 8 "町立静内病院"
Set-PlainValue $ws 1469 9 1
Set-TextValue $ws 1469 10 "2021/8/18"
Set-EmptyTextCell $ws 1469 11
Set-PlainValue $ws 1469 12 1
Set-PlainValue $ws 1469 13 1
Set-PlainValue $ws 1469 14 0
Set-PlainValue $ws 1469 15 0
Set-EmptyTextCell $ws 1469 16
Set-PlainValue $ws 1469 17 "その他"
Set-PlainValue $ws 1469 18 "医療福祉介護"
Set-PlainValue $ws 1469 19 5
# row 1470
Set-PlainValue $ws 1470 1 1469
Set-PlainValue $ws 1470 2 43.19035206
Set-PlainValue $ws 1470 3 140.9941424
Set-PlainValue $ws 1470 4 "北海道"
Set-PlainValue $ws 1470 5 "小樽市"
Set-PlainValue $ws 1470 6 "事業所"
Set-PlainValue $ws 1470 7 3
Set-PlainValue $ws 1470 8 "小樽市役所 財政部"
Set-PlainValue $ws 1470 9 1
Set-TextValue $ws 1470 10 "2021/8/18"
Set-EmptyTextCell $ws 1470 11
Set-PlainValue $ws 1470 12 1
Set-PlainValue $ws 1470 13 0
Set-PlainValue $ws 1470 14 1
Set-PlainValue $ws 1470 15 0
Set-EmptyTextCell $ws 1470 16
Set-PlainValue $ws 1470 17 "その他"
Set-PlainValue $ws 1470 18 "その他"
Set-PlainValue $ws 1470 19 3
# row 1471
Set-PlainValue $ws 1471 1 1470
Set-PlainValue $ws 1471 2 42.33156761
Set-PlainValue $ws 1471 3 140.9527883
Set-PlainValue $ws 1471 4 "北海道"
Set-PlainValue $ws 1471 5 "室蘭市"
Set-PlainValue $ws 1471 6 "小学校"
Set-PlainValue $ws 1471 7 6
Set-PlainValue $ws 1471 8 "室蘭市立みなと小学校"
Set-PlainValue $ws 1471 9 1
Set-TextValue $ws 1471 10 "2021/8/18"
Set-EmptyTextCell $ws 1471 11
Set-PlainValue $ws 1471 12 1
Set-PlainValue $ws 1471 13 0
Set-PlainValue $ws 1471 14 1
Set-PlainValue $ws 1471 15 0
Set-EmptyTextCell $ws 1471 16
Set-PlainValue $ws 1471 17 "その他"
Set-PlainValue $ws 1471 18 "教育/保育"
Set-PlainValue $ws 1471 19 6
# row 1472
Set-PlainValue $ws 1472 1 1471
Set-PlainValue $ws 1472 2 42.37892375
Set-PlainValue $ws 1472 3 141.0345676
Set-PlainValue $ws 1472 4 "北海道"
Set-PlainValue $ws 1472 5 "室蘭市"
Set-PlainValue $ws 1472 6 "大学"
Set-PlainValue $ws 1472 7 9
Set-PlainValue $ws 1472 8 "室蘭工業大学"
Set-PlainValue $ws 1472 9 1
Set-TextValue $ws 1472 10 "2021/8/18"
Set-EmptyTextCell $ws 1472 11
Set-PlainValue $ws 1472 12 1
Set-PlainValue $ws 1472 13 0
Set-PlainValue $ws 1472 14 1
Set-PlainValue $ws 1472 15 0
Set-EmptyTextCell $ws 1472 16
Set-PlainValue $ws 1472 17 "教室"
Set-PlainValue $ws 1472 18 "教育/保育"
Set-PlainValue $ws 1472 19 9
